$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New translation entry: CONTACT.INPUT_SHOP ("Magasin" / "Shop")
$ws.Range("A34").Value = "CONTACT.INPUT_SHOP"
$ws.Range("B34").Value = "Magasin"
$ws.Range("C34").Value = "Shop"

# Update the view state to match where the user was working
$null = $ws.Range("B37").Select()
